$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-23 06:48:33'
$ws.Range('K2').Value = '-0.1 MJ/m2'
$ws.Range('N2').Value = '1.5 °C 6:12 TU'
$ws.Range('O2').Value = '3.3 °C'
$ws.Range('E3').Value = '2026-02-23 06:48:35'
$ws.Range('H3').NumberFormat = '@'
$ws.Range('H3').Value = '41%'
$ws.Range('L3').Value = '32.8 km/h - 223º 6:08 TU'
$ws.Range('E4').Value = '2026-02-23 06:48:38'
$ws.Range('E5').Value = '2026-02-23 06:48:40'
$ws.Range('H5').NumberFormat = '@'
$ws.Range('H5').Value = '36%'
$ws.Range('E6').Value = '2026-02-23 06:48:43'
$ws.Range('N6').Value = '7.7 °C 6:27 TU'
$ws.Range('E7').Value = '2026-02-23 06:48:46'
$ws.Range('K7').Value = '-0.1 MJ/m2'
$ws.Range('N7').Value = '11.0 °C 6:10 TU'
$ws.Range('E8').Value = '2026-02-23 06:48:48'
$ws.Range('O8').Value = '12.8 °C'
$ws.Range('E9').Value = '2026-02-23 06:48:51'
$ws.Range('N9').Value = '4.0 °C 6:23 TU'
$ws.Range('O9').Value = '6.8 °C'
$ws.Range('E10').Value = '2026-02-23 06:48:53'
$ws.Range('E11').Value = '2026-02-23 06:48:56'
$ws.Range('N11').Value = '1.4 °C 6:16 TU'
$ws.Range('O11').Value = '2.7 °C'
$ws.Range('E12').Value = '2026-02-23 06:48:58'
$ws.Range('N12').Value = '3.1 °C 6:19 TU'
$ws.Range('O12').Value = '5.3 °C'
$ws.Range('E13').Value = '2026-02-23 06:49:01'
$ws.Range('L13').Value = '10.8 km/h - 139º 6:10 TU'
$ws.Range('N13').Value = '-3.6 °C 6:29 TU'
$ws.Range('O13').Value = '-1.3 °C'
$ws.Range('E14').Value = '2026-02-23 06:49:04'
$ws.Range('L14').Value = '24.5 km/h - 305º 6:03 TU'
$ws.Range('M14').Value = '9.7 °C 6:08 TU'
$ws.Range('E15').Value = '2026-02-23 06:49:06'
$ws.Range('N15').Value = '4.6 °C 6:13 TU'
$ws.Range('O15').Value = '6.5 °C'
$ws.Range('E16').Value = '2026-02-23 06:49:09'
$ws.Range('E17').Value = '2026-02-23 06:49:11'
$ws.Range('E18').Value = '2026-02-23 06:49:14'
$ws.Range('N18').Value = '1.3 °C 6:24 TU'
$ws.Range('O18').Value = '2.8 °C'
$ws.Range('E19').Value = '2026-02-23 06:49:17'
$ws.Range('N19').Value = '8.1 °C 6:15 TU'
$ws.Range('O19').Value = '9.7 °C'
$ws.Range('E20').Value = '2026-02-23 06:49:19'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '36%'
$ws.Range('E21').Value = '2026-02-23 06:49:22'
$ws.Range('N21').Value = '1.1 °C 6:27 TU'
$ws.Range('O21').Value = '3.3 °C'
$ws.Range('E22').Value = '2026-02-23 06:49:24'
$ws.Range('H22').NumberFormat = '@'
$ws.Range('H22').Value = '24%'
$ws.Range('E23').Value = '2026-02-23 06:49:27'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '25%'
$ws.Range('L23').Value = '36.7 km/h - 336º 6:24 TU'
$ws.Range('E24').Value = '2026-02-23 06:49:30'
$ws.Range('N24').Value = '0.1 °C 6:16 TU'
$ws.Range('O24').Value = '2.1 °C'
$ws.Range('E25').Value = '2026-02-23 06:49:32'
$ws.Range('E26').Value = '2026-02-23 06:49:35'
$ws.Range('O26').Value = '6.4 °C'
$ws.Range('E27').Value = '2026-02-23 06:49:37'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '33%'
$ws.Range('O27').Value = '3.5 °C'
$ws.Range('E28').Value = '2026-02-23 06:49:40'
$ws.Range('N28').Value = '1.3 °C 6:28 TU'
$ws.Range('O28').Value = '3.5 °C'
$ws.Range('E29').Value = '2026-02-23 06:49:42'
$ws.Range('O29').Value = '4.1 °C'
$ws.Range('E30').Value = '2026-02-23 06:49:45'
$ws.Range('J30').Value = '1025.5 hPa'
$ws.Range('L30').Value = '18.4 km/h - 41º 6:29 TU'
$ws.Range('E31').Value = '2026-02-23 06:49:48'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '48%'
$ws.Range('E32').Value = '2026-02-23 06:49:50'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '96%'
$ws.Range('O32').Value = '1.4 °C'
$ws.Range('E33').Value = '2026-02-23 06:49:53'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '67%'
$ws.Range('O33').Value = '2.2 °C'
$ws.Range('E34').Value = '2026-02-23 06:49:56'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '46%'
$ws.Range('L34').Value = '24.1 km/h - 30º 6:26 TU'
$ws.Range('M34').Value = '5.0 °C 6:29 TU'
$ws.Range('O34').Value = '2.1 °C'
$ws.Range('E35').Value = '2026-02-23 06:49:58'
$ws.Range('E36').Value = '2026-02-23 06:50:01'
$ws.Range('J36').Value = '1025.2 hPa'
$ws.Range('O36').Value = '6.7 °C'
$ws.Range('E37').Value = '2026-02-23 06:50:04'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '81%'
$ws.Range('N37').Value = '0.8 °C 6:20 TU'
$ws.Range('O37').Value = '3.5 °C'
$ws.Range('E38').Value = '2026-02-23 06:50:06'
$ws.Range('E39').Value = '2026-02-23 06:50:09'
$ws.Range('O39').Value = '3.6 °C'
$ws.Range('E40').Value = '2026-02-23 06:50:11'
$ws.Range('N40').Value = '-0.1 °C 6:25 TU'
$ws.Range('O40').Value = '1.8 °C'
$ws.Range('E41').Value = '2026-02-23 06:50:14'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '87%'
$ws.Range('J41').Value = '1024.9 hPa'
$ws.Range('N41').Value = '5.4 °C 6:23 TU'
$ws.Range('O41').Value = '7.1 °C'
$ws.Range('E42').Value = '2026-02-23 06:50:17'
$ws.Range('E43').Value = '2026-02-23 06:50:19'
$ws.Range('N43').Value = '1.5 °C 6:18 TU'
$ws.Range('O43').Value = '3.6 °C'
$ws.Range('E44').Value = '2026-02-23 06:50:22'
$ws.Range('K44').Value = '-0.1 MJ/m2'
$ws.Range('O44').Value = '2.5 °C'
$ws.Range('E45').Value = '2026-02-23 06:50:24'
$ws.Range('E46').Value = '2026-02-23 06:50:27'
$ws.Range('N46').Value = '0.3 °C 6:09 TU'
$ws.Range('O46').Value = '1.9 °C'
